$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 6000
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").Value = $null

$ws.Range("H70").Value = 1997.5
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 1997.5
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 5992.5
$ws.Range("N70").Value = -6532.5
$ws.Range("M70").Value = $null

$ws.Range("H73").Value = 1997.5
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 1997.5
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 5992.5
$ws.Range("N73").Value = -7864.5
$ws.Range("M73").Value = $null

$ws.Range("H80").Value = 444.14285
$ws.Range("I80").Value = 341.9091
$ws.Range("K80").Value = 1025.7273
$ws.Range("M80").Value = -27.72730000000001

$ws.Range("H83").Value = 444.14285
$ws.Range("I83").Value = 341.9091
$ws.Range("K83").Value = 3077.1819
$ws.Range("M83").Value = 1914.8181

$ws.Range("H116").Value = 10148.071
$ws.Range("I116").Value = 13711.667
$ws.Range("J116").Value = 3733.6
$ws.Range("K116").Value = 13711.667
$ws.Range("L116").Value = 3733.6
$ws.Range("M116").Value = -10269.667
$ws.Range("N116").Value = -10617.6

$ws.Range("H124").Value = 50780
$ws.Range("J124").Value = 50780
$ws.Range("L124").Value = 50780
$ws.Range("N124").Value = -60600

$ws.Range("H129").Value = 1181.625
$ws.Range("I129").Value = 603.36365
$ws.Range("J129").Value = 1353.5405
$ws.Range("K129").Value = 1810.09095
$ws.Range("L129").Value = 4060.6215
$ws.Range("M129").Value = 3189.90905
$ws.Range("N129").Value = -14060.6215

$ws.Range("H137").Value = 3944.0356
$ws.Range("I137").Value = 2697.72
$ws.Range("J137").Value = 14330
$ws.Range("K137").Value = 8093.16
$ws.Range("L137").Value = 42990
$ws.Range("M137").Value = -5543.16
$ws.Range("N137").Value = -48090

$ws.Range("H138").Value = 2314.8276
$ws.Range("I138").Value = 2001.3478
$ws.Range("J138").Value = 2520.8286
$ws.Range("K138").Value = 6004.0434
$ws.Range("L138").Value = 7562.485799999999
$ws.Range("M138").Value = -864.0434000000005
$ws.Range("N138").Value = -17842.4858

$ws.Range("H141").Value = 6341.3335
$ws.Range("I141").Value = 2827
$ws.Range("K141").Value = 8481
$ws.Range("M141").Value = -3301

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H30").Value = 18898
$ws.Range("I30").Value = 3622.5
$ws.Range("J30").Value = 80000
$ws.Range("K30").Value = 3622.5
$ws.Range("L30").Value = 80000
$ws.Range("M30").Value = -3472.5
$ws.Range("N30").Value = -80300

$ws.Range("H74").Value = 2998.75
$ws.Range("I74").Value = 2855.7144
$ws.Range("J74").Value = 4000
$ws.Range("K74").Value = 2855.7144
$ws.Range("L74").Value = 4000
$ws.Range("M74").Value = -1981.7144
$ws.Range("N74").Value = -5748

$ws.Range("H77").Value = 2998.75
$ws.Range("I77").Value = 2855.7144
$ws.Range("J77").Value = 4000
$ws.Range("K77").Value = 14278.572
$ws.Range("L77").Value = 20000
$ws.Range("M77").Value = -9910.572
$ws.Range("N77").Value = -28736

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 3422.2
$ws.Range("I8").Value = 3422.2
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 3422.2
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -3282.2
$ws.Range("N8").Value = $null

$ws.Range("H86").Value = 166669580
$ws.Range("I86").Value = 166669580
$ws.Range("K86").Value = 166669580
$ws.Range("M86").Value = -166668457

$ws.Range("H89").Value = 166669580
$ws.Range("I89").Value = 166669580
$ws.Range("K89").Value = 833347900
$ws.Range("M89").Value = -833342284

$ws.Range("H94").Value = 1899.8
$ws.Range("I94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("M94").Value = $null

$ws.Range("H134").Value = 3335.8696
$ws.Range("I134").Value = 3492
$ws.Range("J134").Value = 3192.75
$ws.Range("K134").Value = 10476
$ws.Range("L134").Value = 9578.25
$ws.Range("M134").Value = -7941
$ws.Range("N134").Value = -14648.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6676.048
$ws.Range("I31").Value = 1567.0769
$ws.Range("J31").Value = 8966.275
$ws.Range("K31").Value = 1567.0769
$ws.Range("L31").Value = 8966.275
$ws.Range("M31").Value = -1272.0769
$ws.Range("N31").Value = -9556.275

$ws.Range("H34").Value = 6676.048
$ws.Range("I34").Value = 1567.0769
$ws.Range("J34").Value = 8966.275
$ws.Range("K34").Value = 1567.0769
$ws.Range("L34").Value = 8966.275
$ws.Range("M34").Value = -1365.0769
$ws.Range("N34").Value = -9370.275

$ws.Range("H134").Value = 2087.9524
$ws.Range("I134").Value = 2010.5
$ws.Range("J134").Value = 2242.8572
$ws.Range("K134").Value = 6031.5
$ws.Range("L134").Value = 6728.571599999999
$ws.Range("M134").Value = -3496.5
$ws.Range("N134").Value = -11798.5716

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = $null
$ws.Range("N11").Value = $null

$ws.Range("H19").Value = 3414.6667
$ws.Range("J19").Value = 3597.6
$ws.Range("L19").Value = 10792.8
$ws.Range("N19").Value = -11140.8

$ws.Range("H117").Value = 1174.6
$ws.Range("I117").Value = 256.57144
$ws.Range("K117").Value = 769.71432
$ws.Range("M117").Value = 2672.28568

$ws.Range("H131").Value = 1109.174
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 1109.174
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 3327.522
$ws.Range("N131").Value = -13407.522
$ws.Range("M131").Value = $null

$ws.Range("H139").Value = 4178.0264
$ws.Range("I139").Value = 2566.6667
$ws.Range("J139").Value = 4480.1562
$ws.Range("K139").Value = 7700.000100000001
$ws.Range("L139").Value = 13440.4686
$ws.Range("M139").Value = -2560.000100000001
$ws.Range("N139").Value = -23720.4686

$ws.Range("H140").Value = 2126.7856
$ws.Range("I140").Value = 1932.7142
$ws.Range("K140").Value = 5798.142599999999
$ws.Range("M140").Value = -618.1425999999992

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 533335.3
$ws.Range("J11").Value = 300003
$ws.Range("L11").Value = 300003
$ws.Range("N11").Value = -300281

$ws.Range("I80").Value = 56558604
$ws.Range("J80").Value = 102000
$ws.Range("K80").Value = 56558604
$ws.Range("L80").Value = 102000
$ws.Range("M80").Value = -56557606
$ws.Range("N80").Value = -103996

$ws.Range("I83").Value = 56558604
$ws.Range("J83").Value = 102000
$ws.Range("K83").Value = 282793020
$ws.Range("L83").Value = 510000
$ws.Range("M83").Value = -282788028
$ws.Range("N83").Value = -519984

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2850.7727
$ws.Range("I132").Value = 1993.3846
$ws.Range("J132").Value = 4089.2222
$ws.Range("K132").Value = 5980.1538
$ws.Range("L132").Value = 12267.6666
$ws.Range("M132").Value = -3450.1538
$ws.Range("N132").Value = -17327.6666

$ws.Range("H136").Value = 23812810
$ws.Range("I136").Value = 2600.6
$ws.Range("J136").Value = 83338340
$ws.Range("K136").Value = 7801.799999999999
$ws.Range("L136").Value = 250015020
$ws.Range("M136").Value = -5251.799999999999
$ws.Range("N136").Value = -250020120

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 78679
$ws.Range("J46").Value = 78679
$ws.Range("L46").Value = 78679
$ws.Range("N46").Value = -79141

$ws.Range("H132").Value = 5750031
$ws.Range("I132").Value = 4413.222
$ws.Range("J132").Value = 8335559
$ws.Range("K132").Value = 13239.666
$ws.Range("L132").Value = 25006677
$ws.Range("M132").Value = -10709.666
$ws.Range("N132").Value = -25011737

$ws.Range("H134").Value = 78679
$ws.Range("J134").Value = 78679
$ws.Range("L134").Value = 236037
$ws.Range("N134").Value = -241107

$ws.Range("H136").Value = 2894.6453
$ws.Range("I136").Value = 2868.8948
$ws.Range("J136").Value = 2935.4167
$ws.Range("K136").Value = 8606.6844
$ws.Range("L136").Value = 8806.250100000001
$ws.Range("M136").Value = -6056.6844
$ws.Range("N136").Value = -13906.2501
